# Saldo.xlsx update:
#  - Remove the rows for several accounts that dropped out of this export
#    (005607508/CARLOS, 004472404/DILSON, 004870019/MARIA, 004384167/DOUGLAS,
#     004567324/FRANCISCO, 003894173/ANDREA, 004242237/MARIAH)
#  - Update account 004431546 (GABRIELA)'s balance to its new value

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$accountsToRemove = @(
    "005607508",
    "004472404",
    "004870019",
    "004384167",
    "004567324",
    "003894173",
    "004242237"
)

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Walk bottom-to-top so deleting a row never invalidates the row numbers
# still left to inspect above it.
for ($r = $lastRow; $r -ge 1; $r--) {
    $account = [string]$ws.Cells.Item($r, 1).Text

    if ($accountsToRemove -contains $account) {
        $ws.Rows.Item($r).Delete()
    }
    elseif ($account -eq "004431546") {
        $ws.Cells.Item($r, 3).Value = 252082.89
    }
}
